$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range extent (header + data rows, columns A:G)
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Swap columns A and B for every row (1 through lastRow), which swaps
# the header labels ("Date" <-> "EPU_World") as well as all the numeric
# data beneath them.
for ($r = 1; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $valA = $cellA.Value2
    $valB = $cellB.Value2
    $cellA.Value2 = $valB
    $cellB.Value2 = $valA
}
